$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.480748653411865
$ws.Range("B1").Value = 1.871488451957703
$ws.Range("C1").Value = 1.997957706451416
$ws.Range("D1").Value = 1.586897253990173
$ws.Range("E1").Value = 1.389355301856995
